$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Table2" currently spans A1:E2 (header + 1 data row).
# Add a new data row for the "2620. Counter" LeetCode problem.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

# Add the hyperlink for the new row's Link cell first (this is how the link
# was actually entered), then fill in the rest of the row.
$ws.Hyperlinks.Add($ws.Range("E3"), "https://leetcode.com/problems/counter/solutions/3491300/day2-o-1-understanding-closure-in-easy-way-and-its-practical-uses/?envType=study-plan-v2&envId=30-days-of-javascript") | Out-Null

$ws.Range("A3").Value = "2620. Counter"
$ws.Range("B3").Value = "Easy"
$ws.Range("C3").Value = "Closures"
$ws.Range("D3").Value = "Use closure"

# Match the formatting of the corresponding cells on the row above.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Restore the last active selection.
$ws.Range("D13").Select() | Out-Null
